$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on Price (D) cells whose new values would otherwise be
# auto-converted to numbers by Excel, so they remain text just like the source data.
$textCells = @("D4", "D5", "D6", "D9", "D10", "D13", "D18", "D19", "D20", "D21", "D23", "D28", "D30", "D31", "D34", "D35", "D36", "D37", "D38", "D41", "D42", "D43", "D44", "D48", "D49", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply updated cell values
$ws.Range("D2").Value = "62.934.04"
$ws.Range("E2").Value = "  +2.23%  "
$ws.Range("D3").Value = "3.481.27"
$ws.Range("E3").Value = "  +2.70%  "
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "582.20"
$ws.Range("E5").Value = "  +1.08%  "
$ws.Range("D6").Value = "147.74"
$ws.Range("E6").Value = "  +4.72%  "
$ws.Range("E7").Value = "  -0.07%  "
$ws.Range("E8").Value = "  +1.37%  "
$ws.Range("D9").Value = "7.65"
$ws.Range("E9").Value = "  -1.16%  "
$ws.Range("D10").Value = "0.126"
$ws.Range("E10").Value = "  +2.46%  "
$ws.Range("E11").Value = "  +3.83%  "
$ws.Range("D12").Value = "4.077.67"
$ws.Range("E12").Value = "  +2.73%  "
$ws.Range("D13").Value = "29.87"
$ws.Range("E13").Value = "  +5.00%  "
$ws.Range("E14").Value = "  +0.53%  "
$ws.Range("D15").Value = "3.479.43"
$ws.Range("E15").Value = "  +2.86%  "
$ws.Range("E16").Value = "  +1.07%  "
$ws.Range("D17").Value = "62.950.94"
$ws.Range("E17").Value = "  +2.33%  "
$ws.Range("D18").Value = "6.36"
$ws.Range("E18").Value = "  +3.77%  "
$ws.Range("D19").Value = "14.45"
$ws.Range("E19").Value = "  +5.77%  "
$ws.Range("D20").Value = "9.35"
$ws.Range("E20").Value = "  +4.18%  "
$ws.Range("D21").Value = "390.86"
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("E22").Value = "  +2.68%  "
$ws.Range("D23").Value = "75.09"
$ws.Range("E23").Value = "  -0.30%  "
$ws.Range("E24").Value = "  +0.05%  "
$ws.Range("D25").Value = "3.622.82"
$ws.Range("E25").Value = "  +2.58%  "
$ws.Range("E26").Value = "  +3.14%  "
$ws.Range("E27").Value = "  -6.73%  "
$ws.Range("D28").Value = "7.70"
$ws.Range("E28").Value = "  +6.01%  "
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").Value = "8.27"
$ws.Range("E30").Value = "  +2.92%  "
$ws.Range("D31").Value = "2.16"
$ws.Range("E31").Value = "  +0.65%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").Value = "23.85"
$ws.Range("E34").Value = "  +2.52%  "
$ws.Range("B35").Value = "NEARProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D35").Value = "5.32"
$ws.Range("E35").Value = "  +5.23%  "
$ws.Range("B36").Value = "Aptos"
$ws.Range("C36").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D36").Value = "7.13"
$ws.Range("E36").Value = "  +2.90%  "
$ws.Range("D37").Value = "31.78"
$ws.Range("E37").Value = "  +22.23%  "
$ws.Range("D38").Value = "171.34"
$ws.Range("E38").Value = "  +1.80%  "
$ws.Range("E39").Value = "  +7.15%  "
$ws.Range("D40").Value = "3.519.33"
$ws.Range("E40").Value = "  +2.76%  "
$ws.Range("D41").Value = "0.0772"
$ws.Range("E41").Value = "  +0.23%  "
$ws.Range("D42").Value = "0.804"
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "4.53"
$ws.Range("E43").Value = "  +2.40%  "
$ws.Range("B44").Value = "OKB"
$ws.Range("C44").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D44").Value = "42.30"
$ws.Range("E44").Value = "  -0.46%  "
$ws.Range("E45").Value = "  +3.79%  "
$ws.Range("E46").Value = "  +4.31%  "
$ws.Range("D47").Value = "2.602.74"
$ws.Range("E47").Value = "  +5.77%  "
$ws.Range("D48").Value = "23.70"
$ws.Range("E48").Value = "  +3.24%  "
$ws.Range("D49").Value = "2.28"
$ws.Range("E49").Value = "  +11.10%  "
$ws.Range("E50").Value = "  +1.97%  "
$ws.Range("D51").Value = "0.0270"
$ws.Range("E51").Value = "  +2.61%  "
